# Atualização de bases das ligas, do dia: 26-02-2024 às 22:04
#
# For each listed row pair, swap the contents of column B and columns F:AC
# (everything except the id/Div/Div Original Name/Date columns A, C, D, E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns to swap for every pair: B (2), then F..AC (6..29)
$cols = @(2) + @(6..29)

# Row pairs (1-based worksheet row numbers) whose B,F:AC data must be swapped
$rowPairs = @(
    @(28, 29),
    @(47, 48),
    @(55, 56),
    @(131, 132),
    @(149, 150),
    @(215, 216),
    @(221, 222),
    @(230, 231)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    foreach ($c in $cols) {
        $cell1 = $ws.Cells.Item($r1, $c)
        $cell2 = $ws.Cells.Item($r2, $c)

        $v1 = $cell1.Value2
        $v2 = $cell2.Value2

        $cell1.Value2 = $v2
        $cell2.Value2 = $v1
    }
}
